# Scheduled-runner data refresh: updates computed Leve profitability figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across several
# job sheets. Source values come from a re-run of the market-data snapshot;
# only raw numbers change here (no formulas, no formatting).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 736.44446
$ws.Range("I15").Value = 736.44446
$ws.Range("K15").Value = 2209.33338
$ws.Range("M15").Value = -2040.33338

# Row 80
$ws.Range("H80").Value = 3750
$ws.Range("J80").Value = 3333.3333
$ws.Range("L80").Value = 9999.999899999999
$ws.Range("N80").Value = -11995.9999

# Row 83
$ws.Range("H83").Value = 3750
$ws.Range("J83").Value = 3333.3333
$ws.Range("L83").Value = 29999.9997
$ws.Range("N83").Value = -39983.9997

# Row 95
$ws.Range("H95").Value = 39333
$ws.Range("J95").Value = 39333
$ws.Range("L95").Value = 39333
$ws.Range("N95").Value = -44825

# Row 98
$ws.Range("H98").Value = 448
$ws.Range("I98").Value = 273
$ws.Range("K98").Value = 273
$ws.Range("M98").Value = 1225

# Row 107
$ws.Range("H107").Value = 541.375
$ws.Range("I107").Value = 565.8461
$ws.Range("K107").Value = 565.8461
$ws.Range("M107").Value = 1354.1539

# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = $null

# Row 122
$ws.Range("H122").Value = 448
$ws.Range("I122").Value = 273
$ws.Range("K122").Value = 819
$ws.Range("M122").Value = 1631

# Row 132
$ws.Range("H132").Value = 860.6667
$ws.Range("I132").Value = 860.6667
$ws.Range("K132").Value = 2582.0001
$ws.Range("M132").Value = -52.0001000000002

# Row 135
$ws.Range("H135").Value = 431.85715
$ws.Range("I135").Value = 435.8
$ws.Range("K135").Value = 3922.2
$ws.Range("M135").Value = -1387.2

# Row 138
$ws.Range("H138").Value = 2386.647
$ws.Range("I138").Value = 2042.3077
$ws.Range("K138").Value = 6126.9231
$ws.Range("M138").Value = -986.9231

$ws = $wb.Worksheets.Item("ARM")
# Row 60
$ws.Range("H60").Value = 9025.5
$ws.Range("I60").Value = 9025.5
$ws.Range("K60").Value = 9025.5
$ws.Range("M60").Value = -8292.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1540.4286
$ws.Range("I31").Value = 1540.4286
$ws.Range("K31").Value = 1540.4286
$ws.Range("M31").Value = -1245.4286

# Row 34
$ws.Range("H34").Value = 1540.4286
$ws.Range("I34").Value = 1540.4286
$ws.Range("K34").Value = 1540.4286
$ws.Range("M34").Value = -1338.4286

# Row 95
$ws.Range("H95").Value = 46599.8
$ws.Range("J95").Value = 46599.8
$ws.Range("L95").Value = 46599.8
$ws.Range("N95").Value = -52091.8

# Row 99
$ws.Range("H99").Value = 2001540
$ws.Range("J99").Value = 2500750
$ws.Range("L99").Value = 2500750
$ws.Range("N99").Value = -2503746

# Row 126
$ws.Range("H126").Value = 2001540
$ws.Range("J126").Value = 2500750
$ws.Range("L126").Value = 7502250
$ws.Range("N126").Value = -7507190

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 1421
$ws.Range("I46").Value = 1421
$ws.Range("K46").Value = 1421
$ws.Range("M46").Value = -1265

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null

# Row 132
$ws.Range("H132").Value = 1780.6364
$ws.Range("I132").Value = 1590.2
$ws.Range("J132").Value = 3685
$ws.Range("K132").Value = 4770.6
$ws.Range("L132").Value = 11055
$ws.Range("M132").Value = -2240.6
$ws.Range("N132").Value = -16115

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4950.0835
$ws.Range("I7").Value = 4950.0835
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4950.0835
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4838.0835
$ws.Range("N7").Value = $null

# Row 22
$ws.Range("H22").Value = 776.1667
$ws.Range("I22").Value = 708.4
$ws.Range("J22").Value = 824.5714
$ws.Range("K22").Value = 708.4
$ws.Range("L22").Value = 824.5714
$ws.Range("M22").Value = -413.4
$ws.Range("N22").Value = -1414.5714

# Row 27
$ws.Range("H27").Value = 776.1667
$ws.Range("I27").Value = 708.4
$ws.Range("J27").Value = 824.5714
$ws.Range("K27").Value = 708.4
$ws.Range("L27").Value = 824.5714
$ws.Range("M27").Value = -601.4
$ws.Range("N27").Value = -1038.5714

# Row 61
$ws.Range("H61").Value = 2017
$ws.Range("I61").Value = 2034.5
$ws.Range("J61").Value = 1999.5
$ws.Range("K61").Value = 2034.5
$ws.Range("L61").Value = 1999.5
$ws.Range("M61").Value = -1832.5
$ws.Range("N61").Value = -2403.5

# Row 101
$ws.Range("H101").Value = 21687
$ws.Range("J101").Value = 21687
$ws.Range("L101").Value = 21687
$ws.Range("N101").Value = -28177

# Row 113
$ws.Range("H113").Value = 2017
$ws.Range("I113").Value = 2034.5
$ws.Range("J113").Value = 1999.5
$ws.Range("K113").Value = 2034.5
$ws.Range("L113").Value = 1999.5
$ws.Range("M113").Value = 135.5
$ws.Range("N113").Value = -6339.5

# Row 126
$ws.Range("H126").Value = 4950.0835
$ws.Range("I126").Value = 4950.0835
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14850.2505
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12380.2505
$ws.Range("N126").Value = $null

# Row 132
$ws.Range("H132").Value = 3240.7778
$ws.Range("I132").Value = 3026.7144
$ws.Range("J132").Value = 3990
$ws.Range("K132").Value = 9080.143199999999
$ws.Range("L132").Value = 11970
$ws.Range("M132").Value = -6550.143199999999
$ws.Range("N132").Value = -17030

# Row 136
$ws.Range("H136").Value = 1320.8
$ws.Range("J136").Value = 1400
$ws.Range("L136").Value = 4200
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 27500
$ws.Range("J101").Value = 27500
$ws.Range("L101").Value = 27500
$ws.Range("N101").Value = -33990

# Row 107
$ws.Range("H107").Value = 1171.625
$ws.Range("I107").Value = 894.6667
$ws.Range("K107").Value = 2684.0001
$ws.Range("M107").Value = -764.0001000000002
